$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-07-14 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-15 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("76-53=", $true, $false, $false, $false, $false, $true, 1, $false, "59-47=", 2) | Out-Null
$d.Content.Find.Execute("33+1=", $true, $false, $false, $false, $false, $true, 1, $false, "58-30=", 2) | Out-Null
$d.Content.Find.Execute("22+30=", $true, $false, $false, $false, $false, $true, 1, $false, "7+19=", 2) | Out-Null
$d.Content.Find.Execute("44+30=", $true, $false, $false, $false, $false, $true, 1, $false, "60-48=", 2) | Out-Null
$d.Content.Find.Execute("22-5=", $true, $false, $false, $false, $false, $true, 1, $false, "76+20=", 2) | Out-Null
$d.Content.Find.Execute("69+11=", $true, $false, $false, $false, $false, $true, 1, $false, "23+10=", 2) | Out-Null
$d.Content.Find.Execute("17+19=", $true, $false, $false, $false, $false, $true, 1, $false, "62-42=", 2) | Out-Null
$d.Content.Find.Execute("64+3=", $true, $false, $false, $false, $false, $true, 1, $false, "76-57=", 2) | Out-Null
$d.Content.Find.Execute("17+37=", $true, $false, $false, $false, $false, $true, 1, $false, "21+30=", 2) | Out-Null
$d.Content.Find.Execute("67+10=", $true, $false, $false, $false, $false, $true, 1, $false, "21+20=", 2) | Out-Null
$d.Content.Find.Execute("6+31=", $true, $false, $false, $false, $false, $true, 1, $false, "49+38=", 2) | Out-Null
$d.Content.Find.Execute("2+33=", $true, $false, $false, $false, $false, $true, 1, $false, "63-60=", 2) | Out-Null
$d.Content.Find.Execute("37+27=", $true, $false, $false, $false, $false, $true, 1, $false, "59+37=", 2) | Out-Null
$d.Content.Find.Execute("79+15=", $true, $false, $false, $false, $false, $true, 1, $false, "90-80=", 2) | Out-Null
$d.Content.Find.Execute("2+58=", $true, $false, $false, $false, $false, $true, 1, $false, "58+7=", 2) | Out-Null
$d.Content.Find.Execute("86-15=", $true, $false, $false, $false, $false, $true, 1, $false, "55+19=", 2) | Out-Null
$d.Content.Find.Execute("72+24=", $true, $false, $false, $false, $false, $true, 1, $false, "47-38=", 2) | Out-Null
$d.Content.Find.Execute("10+43=", $true, $false, $false, $false, $false, $true, 1, $false, "32+33=", 2) | Out-Null
$d.Content.Find.Execute("12+12=", $true, $false, $false, $false, $false, $true, 1, $false, "22+66=", 2) | Out-Null
$d.Content.Find.Execute("94-62=", $true, $false, $false, $false, $false, $true, 1, $false, "68-27=", 2) | Out-Null
$d.Content.Find.Execute("30+10=", $true, $false, $false, $false, $false, $true, 1, $false, "50+26=", 2) | Out-Null
$d.Content.Find.Execute("16+62=", $true, $false, $false, $false, $false, $true, 1, $false, "42-3=", 2) | Out-Null
$d.Content.Find.Execute("96-46=", $true, $false, $false, $false, $false, $true, 1, $false, "1+93=", 2) | Out-Null
$d.Content.Find.Execute("51-34=", $true, $false, $false, $false, $false, $true, 1, $false, "4+27=", 2) | Out-Null
$d.Content.Find.Execute("81-52=", $true, $false, $false, $false, $false, $true, 1, $false, "24+63=", 2) | Out-Null
$d.Content.Find.Execute("14+3=", $true, $false, $false, $false, $false, $true, 1, $false, "40-37=", 2) | Out-Null
$d.Content.Find.Execute("74-18=", $true, $false, $false, $false, $false, $true, 1, $false, "92-29=", 2) | Out-Null
$d.Content.Find.Execute("91-57=", $true, $false, $false, $false, $false, $true, 1, $false, "0+31=", 2) | Out-Null
$d.Content.Find.Execute("72+15=", $true, $false, $false, $false, $false, $true, 1, $false, "10+70=", 2) | Out-Null
$d.Content.Find.Execute("29+43=", $true, $false, $false, $false, $false, $true, 1, $false, "43-22=", 2) | Out-Null
$d.Content.Find.Execute("44+27=", $true, $false, $false, $false, $false, $true, 1, $false, "9-7=", 2) | Out-Null
$d.Content.Find.Execute("11+56=", $true, $false, $false, $false, $false, $true, 1, $false, "34-31=", 2) | Out-Null
$d.Content.Find.Execute("67-3=", $true, $false, $false, $false, $false, $true, 1, $false, "62+31=", 2) | Out-Null
$d.Content.Find.Execute("16+64=", $true, $false, $false, $false, $false, $true, 1, $false, "75+0=", 2) | Out-Null
$d.Content.Find.Execute("33+31=", $true, $false, $false, $false, $false, $true, 1, $false, "38+48=", 2) | Out-Null
$d.Content.Find.Execute("78+7=", $true, $false, $false, $false, $false, $true, 1, $false, "60+36=", 2) | Out-Null
$d.Content.Find.Execute("30+22=", $true, $false, $false, $false, $false, $true, 1, $false, "58-23=", 2) | Out-Null
$d.Content.Find.Execute("81+11=", $true, $false, $false, $false, $false, $true, 1, $false, "54+11=", 2) | Out-Null
$d.Content.Find.Execute("68-24=", $true, $false, $false, $false, $false, $true, 1, $false, "51-18=", 2) | Out-Null
$d.Content.Find.Execute("54+3=", $true, $false, $false, $false, $false, $true, 1, $false, "55-14=", 2) | Out-Null
$d.Content.Find.Execute("11-3=", $true, $false, $false, $false, $false, $true, 1, $false, "31-28=", 2) | Out-Null
$d.Content.Find.Execute("71-25=", $true, $false, $false, $false, $false, $true, 1, $false, "35-9=", 2) | Out-Null
$d.Content.Find.Execute("60+10=", $true, $false, $false, $false, $false, $true, 1, $false, "79-62=", 2) | Out-Null
$d.Content.Find.Execute("33+6=", $true, $false, $false, $false, $false, $true, 1, $false, "13+49=", 2) | Out-Null
$d.Content.Find.Execute("93-80=", $true, $false, $false, $false, $false, $true, 1, $false, "91-9=", 2) | Out-Null
$d.Content.Find.Execute("34-4=", $true, $false, $false, $false, $false, $true, 1, $false, "54+16=", 2) | Out-Null
$d.Content.Find.Execute("5+51=", $true, $false, $false, $false, $false, $true, 1, $false, "28+3=", 2) | Out-Null
$d.Content.Find.Execute("4+53=", $true, $false, $false, $false, $false, $true, 1, $false, "95-50=", 2) | Out-Null
$d.Content.Find.Execute("12-5=", $true, $false, $false, $false, $false, $true, 1, $false, "65-44=", 2) | Out-Null
$d.Content.Find.Execute("52+15=", $true, $false, $false, $false, $false, $true, 1, $false, "58+19=", 2) | Out-Null
$d.Content.Find.Execute("88-37=", $true, $false, $false, $false, $false, $true, 1, $false, "25-25=", 2) | Out-Null
$d.Content.Find.Execute("9+82=", $true, $false, $false, $false, $false, $true, 1, $false, "34+0=", 2) | Out-Null
$d.Content.Find.Execute("42+17=", $true, $false, $false, $false, $false, $true, 1, $false, "17+50=", 2) | Out-Null
$d.Content.Find.Execute("64-39=", $true, $false, $false, $false, $false, $true, 1, $false, "14+20=", 2) | Out-Null
$d.Content.Find.Execute("0+32=", $true, $false, $false, $false, $false, $true, 1, $false, "83-78=", 2) | Out-Null
$d.Content.Find.Execute("99-54=", $true, $false, $false, $false, $false, $true, 1, $false, "20+30=", 2) | Out-Null
$d.Content.Find.Execute("66-2=", $true, $false, $false, $false, $false, $true, 1, $false, "8+83=", 2) | Out-Null
$d.Content.Find.Execute("20+33=", $true, $false, $false, $false, $false, $true, 1, $false, "23+68=", 2) | Out-Null
$d.Content.Find.Execute("35+12=", $true, $false, $false, $false, $false, $true, 1, $false, "96-40=", 2) | Out-Null
$d.Content.Find.Execute("28-3=", $true, $false, $false, $false, $false, $true, 1, $false, "91-89=", 2) | Out-Null
$d.Content.Find.Execute("68-36=", $true, $false, $false, $false, $false, $true, 1, $false, "67-51=", 2) | Out-Null
$d.Content.Find.Execute("24+18=", $true, $false, $false, $false, $false, $true, 1, $false, "37-14=", 2) | Out-Null
$d.Content.Find.Execute("57+14=", $true, $false, $false, $false, $false, $true, 1, $false, "46-0=", 2) | Out-Null
$d.Content.Find.Execute("1+29=", $true, $false, $false, $false, $false, $true, 1, $false, "85-47=", 2) | Out-Null
$d.Content.Find.Execute("49-13=", $true, $false, $false, $false, $false, $true, 1, $false, "67-8=", 2) | Out-Null
$d.Content.Find.Execute("60-18=", $true, $false, $false, $false, $false, $true, 1, $false, "68+16=", 2) | Out-Null
$d.Content.Find.Execute("1+31=", $true, $false, $false, $false, $false, $true, 1, $false, "6+49=", 2) | Out-Null
$d.Content.Find.Execute("66-33=", $true, $false, $false, $false, $false, $true, 1, $false, "38-10=", 2) | Out-Null
$d.Content.Find.Execute("8+31=", $true, $false, $false, $false, $false, $true, 1, $false, "13+20=", 2) | Out-Null
$d.Content.Find.Execute("36+42=", $true, $false, $false, $false, $false, $true, 1, $false, "73+14=", 2) | Out-Null
$d.Content.Find.Execute("5+15=", $true, $false, $false, $false, $false, $true, 1, $false, "67+31=", 2) | Out-Null
$d.Content.Find.Execute("16+33=", $true, $false, $false, $false, $false, $true, 1, $false, "77+2=", 2) | Out-Null
$d.Content.Find.Execute("54-31=", $true, $false, $false, $false, $false, $true, 1, $false, "75-74=", 2) | Out-Null
$d.Content.Find.Execute("58-24=", $true, $false, $false, $false, $false, $true, 1, $false, "73-27=", 2) | Out-Null
$d.Content.Find.Execute("78-76=", $true, $false, $false, $false, $false, $true, 1, $false, "55+15=", 2) | Out-Null
$d.Content.Find.Execute("98-79=", $true, $false, $false, $false, $false, $true, 1, $false, "41+25=", 2) | Out-Null
$d.Content.Find.Execute("26+44=", $true, $false, $false, $false, $false, $true, 1, $false, "2+31=", 2) | Out-Null
$d.Content.Find.Execute("76-14=", $true, $false, $false, $false, $false, $true, 1, $false, "63-11=", 2) | Out-Null
$d.Content.Find.Execute("20+61=", $true, $false, $false, $false, $false, $true, 1, $false, "60-39=", 2) | Out-Null
$d.Content.Find.Execute("7+46=", $true, $false, $false, $false, $false, $true, 1, $false, "61+4=", 2) | Out-Null
$d.Content.Find.Execute("26+53=", $true, $false, $false, $false, $false, $true, 1, $false, "68+14=", 2) | Out-Null
$d.Content.Find.Execute("21+16=", $true, $false, $false, $false, $false, $true, 1, $false, "27+19=", 2) | Out-Null
$d.Content.Find.Execute("99-44=", $true, $false, $false, $false, $false, $true, 1, $false, "50+4=", 2) | Out-Null
$d.Content.Find.Execute("58-32=", $true, $false, $false, $false, $false, $true, 1, $false, "45-25=", 2) | Out-Null
$d.Content.Find.Execute("12+25=", $true, $false, $false, $false, $false, $true, 1, $false, "90-73=", 2) | Out-Null
$d.Content.Find.Execute("8+88=", $true, $false, $false, $false, $false, $true, 1, $false, "86-32=", 2) | Out-Null
$d.Content.Find.Execute("71-68=", $true, $false, $false, $false, $false, $true, 1, $false, "67+24=", 2) | Out-Null
$d.Content.Find.Execute("28+60=", $true, $false, $false, $false, $false, $true, 1, $false, "22+47=", 2) | Out-Null
$d.Content.Find.Execute("82-45=", $true, $false, $false, $false, $false, $true, 1, $false, "37+14=", 2) | Out-Null
$d.Content.Find.Execute("12+21=", $true, $false, $false, $false, $false, $true, 1, $false, "96-6=", 2) | Out-Null
$d.Content.Find.Execute("98-52=", $true, $false, $false, $false, $false, $true, 1, $false, "84-31=", 2) | Out-Null
$d.Content.Find.Execute("23+60=", $true, $false, $false, $false, $false, $true, 1, $false, "75-64=", 2) | Out-Null
$d.Content.Find.Execute("35+8=", $true, $false, $false, $false, $false, $true, 1, $false, "68+2=", 2) | Out-Null
$d.Content.Find.Execute("98-54=", $true, $false, $false, $false, $false, $true, 1, $false, "75-70=", 2) | Out-Null
$d.Content.Find.Execute("74-42=", $true, $false, $false, $false, $false, $true, 1, $false, "35+13=", 2) | Out-Null
$d.Content.Find.Execute("46-4=", $true, $false, $false, $false, $false, $true, 1, $false, "54+31=", 2) | Out-Null
$d.Content.Find.Execute("40+30=", $true, $false, $false, $false, $false, $true, 1, $false, "82-56=", 2) | Out-Null
$d.Content.Find.Execute("70-60=", $true, $false, $false, $false, $false, $true, 1, $false, "8+82=", 2) | Out-Null
$d.Content.Find.Execute("58+33=", $true, $false, $false, $false, $false, $true, 1, $false, "9+19=", 2) | Out-Null
$d.Content.Find.Execute("6+5=", $true, $false, $false, $false, $false, $true, 1, $false, "49-19=", 2) | Out-Null
